$d = $word.ActiveDocument

# --- Change 1: after the "Chức năng 3" paragraph's sentence (ending "...qua
#     địa chỉ email."), append a new, separate, non-bold Arial run " - " ---
$r1 = $d.Content
$r1.Find.Execute(
    "qua địa chỉ email.", $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0
) | Out-Null
$r1.Collapse(0)
$pos1 = $r1.Start
$r1.InsertAfter(" - ")
$new1 = $d.Range($pos1, $pos1 + 3)
$new1.Font.Name = "Arial"

# --- Change 2: after "Chức năng 3.1", append a new, separate, bold Arial
#     run " - DONE" ---
$r2 = $d.Content
$r2.Find.Execute(
    "Chức năng 3.1", $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0
) | Out-Null
$r2.Collapse(0)
$pos2 = $r2.Start
$r2.InsertAfter(" - DONE")
$new2 = $d.Range($pos2, $pos2 + 7)
$new2.Font.Name = "Arial"
$new2.Font.Bold = $true

# --- Change 3: move the "_GoBack" bookmark from the end of the "Chức năng
#     3.2" paragraph into the middle of "phút trc 23h" (splitting it into
#     "phút t" | bookmark | "rc 23h") ---
$r3 = $d.Content
$r3.Find.Execute(
    "phút trc 23h", $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$splitPos = $r3.Start + 6

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
